# Add a "unit" column (column D) to the WiscSIMSColumnDictionary sheet,
# populating it with the unit for each dictionary entry, and add several
# new dictionary rows (24-34) that also carry a unit value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D values for the already-existing rows (1-23) ---------------
$ws.Range("D1").Value  = "unit"
$ws.Range("D4").Value  = "permille_VSMOW"
$ws.Range("D5").Value  = "permille"
$ws.Range("D6").Value  = "permille"
$ws.Range("D7").Value  = "permille"
$ws.Range("D8").Value  = "permille"
$ws.Range("D9").Value  = "Gcps"
$ws.Range("D10").Value = "nA"
$ws.Range("D11").Value = "Gcps/nA"
$ws.Range("D12").Value = "Date"
$ws.Range("D13").Value = "Time"
$ws.Range("D14").Value = "micrometers"
$ws.Range("D15").Value = "micrometers"
$ws.Range("D16").Value = "bits"
$ws.Range("D17").Value = "bits"
$ws.Range("D18").Value = "AMU"
$ws.Range("D19").Value = "HydRatio"
$ws.Range("D20").Value = "Mcps"
$ws.Range("D21").Value = "permille_VPDB"
$ws.Range("D22").Value = "cps"
$ws.Range("D23").Value = "Mcps"

# --- New dictionary rows (24-34): ColNames (A) + unit (D) ---------------
$ws.Range("A24").Value = "DATETIME"
$ws.Range("D24").Value = "Date-time"

$ws.Range("A25").Value = "AnalysisLength"
$ws.Range("D25").Value = "minutes"

$ws.Range("A26").Value = "MATERIAL"
$ws.Range("D26").Value = "Sample-std"

$ws.Range("A27").Value = "GROUPNUM"
$ws.Range("D27").Value = "spacedgroup"

$ws.Range("A28").Value = "GUESS.SAMP"
$ws.Range("D28").Value = "ClusterSample"

$ws.Range("A29").Value = "UNIQUEGRP"
$ws.Range("D29").Value = "StdGroup"

$ws.Range("A30").Value = "REL_YIELD"
$ws.Range("D30").Value = "%"

$ws.Range("A31").Value = "REL_Hyd"
$ws.Range("D31").Value = "HydRatio-ratio"

$ws.Range("A32").Value = "BRACKET2SD"
$ws.Range("D32").Value = "permille"

$ws.Range("A33").Value = "STDiso"
$ws.Range("D33").Value = "permille"

$ws.Range("A34").Value = "RegexSTD"
$ws.Range("D34").Value = "RunStd"

# --- Column widths: B narrows, new C gets a width, A stays as-is --------
# (ColumnWidth is in characters; the engine adds ~5/6 of a character when
# converting to the stored xlsx "width" units, so subtract that back out.)
$ws.Columns.Item(2).ColumnWidth = 12 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 20.1640625 - (5/6)

# --- Sheet view: scroll down a bit and move the active selection --------
$ws.Activate()
$ws.Range("F37").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
